# Update StructureDefinition-INT-POS.xlsx to the published CDA FHIR logical
# model with patches #241.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates -------------------------------------------------

# Version
$wsMeta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date
$wsMeta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact
$wsMeta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates (INT_POS.value row, row 4) ----------------------

# Short / Definition / Comments no longer populated for this row
$wsElem.Range("L4").Value = ""
$wsElem.Range("M4").Value = ""
$wsElem.Range("N4").Value = ""

# Base Path stays "INT.value" (unchanged content, renumbering handled by engine)
$wsElem.Range("AF4").Value = "INT.value"

# Constraint(s) no longer carries the ele-1 FHIR constraint text
$wsElem.Range("AJ4").Value = ""
